$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "OK" values from E1:H1
$ws.Range("E1:H1").ClearContents()

# Update A2/B2 remain same values (Cristian / Echevarria) - no change needed there

# Set the merged cell A10:F17 text and style (copy style from C1 which already
# uses the centered alignment style so the same style index is reused)
$ws.Range("C1").Copy()
$ws.Range("A10:F17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A10").Value = "DATOS FUSIONADOS"

# Widen column B to fit the new (longer) content, mirroring the autofit the
# author triggered after entering "DATOS FUSIONADOS" / editing the sheet
$ws.Columns("B").ColumnWidth = 16.45

# Update selection
$ws.Range("I9").Select()
